# Updated cryptos list - refresh Price/Volume(1h) figures and re-sync
# the Quant/TrustWalletToken/PaxDollar row ordering to match the latest
# coinranking.com snapshot.
#
# Note: columns D (Price) and E (Volume(1h)) hold text values that look
# numeric (e.g. "1.000", "0.000007634", "  +0.08%  "). Excel's Range.Value
# setter auto-converts such strings to real numbers/dates, which would
# corrupt the formatting (leading/trailing spaces, trailing zeros, grouped
# "thousand" dots, percent signs, etc.). Prefixing the assigned string with
# a leading apostrophe forces Excel to keep it as literal text, exactly as
# it was stored in the original sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.275.39"
$ws.Range("E2").Value = "'  +0.08%  "
$ws.Range("D3").Value = "'1.869.18"
$ws.Range("E3").Value = "'  +0.27%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'235.08"
$ws.Range("E5").Value = "'  -0.92%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  -0.03%  "
$ws.Range("D7").Value = "'0.4693"
$ws.Range("E7").Value = "'  +0.16%  "
$ws.Range("D8").Value = "'0.2864"
$ws.Range("E8").Value = "'  -0.13%  "
$ws.Range("D9").Value = "'0.06579"
$ws.Range("E9").Value = "'  +0.37%  "
$ws.Range("D10").Value = "'21.78"
$ws.Range("E10").Value = "'  -1.80%  "
$ws.Range("D11").Value = "'0.07991"
$ws.Range("E11").Value = "'  +1.11%  "
$ws.Range("E12").Value = "'  -1.22%  "
$ws.Range("D13").Value = "'1.870.96"
$ws.Range("E13").Value = "'  +0.28%  "
$ws.Range("D14").Value = "'0.6882"
$ws.Range("E14").Value = "'  +0.90%  "
$ws.Range("D15").Value = "'5.114"
$ws.Range("E15").Value = "'  -1.54%  "
$ws.Range("D16").Value = "'269.02"
$ws.Range("E16").Value = "'  -3.41%  "
$ws.Range("D17").Value = "'30.298.53"
$ws.Range("E17").Value = "'  +0.15%  "
$ws.Range("D18").Value = "'14.15"
$ws.Range("E18").Value = "'  +3.80%  "
$ws.Range("D19").Value = "'0.000007634"
$ws.Range("E19").Value = "'  +3.97%  "
$ws.Range("E20").Value = "'  +0.02%  "
$ws.Range("D21").Value = "'2.114.95"
$ws.Range("E21").Value = "'  -0.09%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "'  -0.02%  "
$ws.Range("E23").Value = "'  -2.03%  "
$ws.Range("E24").Value = "'  +0.19%  "
$ws.Range("D25").Value = "'9.394"
$ws.Range("E25").Value = "'  +1.60%  "
$ws.Range("D26").Value = "'167.68"
$ws.Range("E26").Value = "'  -0.21%  "
$ws.Range("D27").Value = "'18.87"
$ws.Range("E27").Value = "'  -1.24%  "
$ws.Range("D28").Value = "'1.948"
$ws.Range("E28").Value = "'  -0.38%  "
$ws.Range("D29").Value = "'1.366"
$ws.Range("E29").Value = "'  -1.17%  "
$ws.Range("D30").Value = "'0.09878"
$ws.Range("E30").Value = "'  +0.20%  "
$ws.Range("D31").Value = "'4.351"
$ws.Range("E31").Value = "'  -0.95%  "
$ws.Range("D32").Value = "'1.456"
$ws.Range("E32").Value = "'  -1.86%  "
$ws.Range("D33").Value = "'4.065"
$ws.Range("E33").Value = "'  -0.26%  "
$ws.Range("D34").Value = "'0.04716"
$ws.Range("E34").Value = "'  -0.84%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("E35").Value = "'  -0.39%  "
$ws.Range("D36").Value = "'0.7015"
$ws.Range("E36").Value = "'  -0.49%  "
$ws.Range("D37").Value = "'2.737"
$ws.Range("E37").Value = "'  +1.10%  "
$ws.Range("D38").Value = "'0.01882"
$ws.Range("E38").Value = "'  +0.00%  "
$ws.Range("D39").Value = "'2.781"
$ws.Range("E39").Value = "'  +5.71%  "
$ws.Range("D40").Value = "'6.260"
$ws.Range("E40").Value = "'  -0.51%  "
$ws.Range("D41").Value = "'72.07"
$ws.Range("E41").Value = "'  -5.13%  "
$ws.Range("D42").Value = "'1.957"
$ws.Range("E42").Value = "'  -0.10%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8422"
$ws.Range("E43").Value = "'  -1.58%  "
$ws.Range("D44").Value = "'0.4174"
$ws.Range("E44").Value = "'  -0.30%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9992"
$ws.Range("E45").Value = "'  -0.06%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'103.12"
$ws.Range("E46").Value = "'  -0.36%  "
$ws.Range("D47").Value = "'7.082"
$ws.Range("E47").Value = "'  -2.10%  "
$ws.Range("D48").Value = "'9.158"
$ws.Range("E48").Value = "'  -0.84%  "
$ws.Range("D49").Value = "'915.71"
$ws.Range("E49").Value = "'  -3.14%  "
$ws.Range("D50").Value = "'34.47"
$ws.Range("E50").Value = "'  +0.49%  "
$ws.Range("D51").Value = "'0.05698"
$ws.Range("E51").Value = "'  +0.89%  "
